$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Domain value on row 2 (was "www.dellservicesevents.com", now "www.dell.com")
$ws.Range("A2").Value = "www.dell.com"

# Update the active selection to H2 (was G2)
$ws.Range("H2").Select()
